$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: fix column layout (col A:B combined -> split, so col B gets its own width) ---
$ws.Columns.Item(2).ColumnWidth = 60.7109375

# --- Step 2: insert a new row at 13 for "Docentes responsaveis" value (label already exists in A12) ---
$ws.Rows.Item(13).Insert()
$ws.Range("A13").Clear()
$ws.Range("B14:C14").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B13").Value = "5840942 - Marco Aurélio Kondracki de Alcântara"
$ws.Range("C13").Value = "5840942 - Marco Aurélio Kondracki de Alcântara"

# --- Step 3: correct the misaligned values that were already shifted by the row insert ---
$ws.Range("B10").Value = "A disciplina tem o objetivo de apresentar ao estudante informações a respeito dos fenômenos químicos que ocorrem nos solos, a fim de capacitá-lo a compreender sua importância na fertilidade do solo e na preservação do ambiente."
$ws.Range("C10").Value = "A disciplina tem o objetivo de apresentar ao estudante informações a respeito dos fenômenos químicos que ocorrem nos solos, a fim de capacitá-lo a compreender sua importância na fertilidade do solo e na preservação do ambiente."

$ws.Range("B14").Value = "Revisão. Composição da fase sólida mineral do solo. Composição da fase sólida orgânica do solo. Solução do solo. Fenômenos de Superfície. Sólidos ácidos e solos afetados por sais. Solos alagados.  Aula prática: Métodos de avaliação da fertilidade do solo. Aula prática: Análises químicas da terra para fins de fertilidade."
$ws.Range("C14").Value = "Revisão. Composição da fase sólida mineral do solo. Composição da fase sólida orgânica do solo. Solução do solo. Fenômenos de Superfície. Sólidos ácidos e solos afetados por sais. Solos alagados.  Aula prática: Métodos de avaliação da fertilidade do solo. Aula prática: Análises químicas da terra para fins de fertilidade."

$ws.Range("B16").Value = "REVISÃO. Conceitos de solo, perfil, composição, características e propriedades. COMPOSIÇÃO DA FASE SÓLIDA MINERAL DO SOLO. Minerais do solo. Principais classes de minerais. Origem das cargas elétricas. COMPOSIÇÃO DA FASE SÓLIDA ORGÂNICA DO SOLO. Composição e estrutura da matéria orgânica do solo, funções e reações, matéria orgânica e sistemas de manejo. SOLUÇÃO DO SOLO. Composição da solução do solo, moléculas orgânicas dissolvidas na solução do solo, concentração e atividade de íons, obtenção da solução do solo. FENÔMENOS DE SUPERFÍCIE. Origem das cargas elétricas, modelos de distribuição das cargas elétricas ao redor das partículas, complexos de superfície, capacidade de troca de cátions e de ânions, adsorção específica, modelos descritivos da adsorção, ponto de carga elétrica zero. Aula prática de campo: Métodos de avaliação da fertilidade do solo: Amostragem de terra: planejamento da amostragem e coletas de amostras de terra. Aula prática de laboratório: Análises químicas da terra para fins de fertilidade: extratores e métodos analíticos."
$ws.Range("C16").Value = "REVISÃO. Conceitos de solo, perfil, composição, características e propriedades. COMPOSIÇÃO DA FASE SÓLIDA MINERAL DO SOLO. Minerais do solo. Principais classes de minerais. Origem das cargas elétricas. COMPOSIÇÃO DA FASE SÓLIDA ORGÂNICA DO SOLO. Composição e estrutura da matéria orgânica do solo, funções e reações, matéria orgânica e sistemas de manejo. SOLUÇÃO DO SOLO. Composição da solução do solo, moléculas orgânicas dissolvidas na solução do solo, concentração e atividade de íons, obtenção da solução do solo. FENÔMENOS DE SUPERFÍCIE. Origem das cargas elétricas, modelos de distribuição das cargas elétricas ao redor das partículas, complexos de superfície, capacidade de troca de cátions e de ânions, adsorção específica, modelos descritivos da adsorção, ponto de carga elétrica zero. Aula prática de campo: Métodos de avaliação da fertilidade do solo: Amostragem de terra: planejamento da amostragem e coletas de amostras de terra. Aula prática de laboratório: Análises químicas da terra para fins de fertilidade: extratores e métodos analíticos."

$ws.Range("B19").Value = "A avaliação será feita mediante duas avaliações escritas de igual peso (P1 e P2). Alternativamente, essas avaliações escritas poderão ser substituídas por duas notas de igual peso (NOTA 1 e NOTA 2). Essas notas serão dadas pela média entre atividades desenvolvidas em aula, trabalhos e relatórios de aulas práticas."
$ws.Range("C19").Value = "A avaliação será feita mediante duas avaliações escritas de igual peso (P1 e P2). Alternativamente, essas avaliações escritas poderão ser substituídas por duas notas de igual peso (NOTA 1 e NOTA 2). Essas notas serão dadas pela média entre atividades desenvolvidas em aula, trabalhos e relatórios de aulas práticas."

$ws.Range("B20").Value = "O aluno poderá optar por dois critérios de avaliação:Critério 1: NF = (P1+P2)/2; ouCritério 2: NF = (NOTA 1 + NOTA 2)/2Sendo P1 e P2 avaliações escritas e NOTA 1 e NOTA 2 obtidas em atividades desenvolvidas em aula, trabalhos e relatórios de aulas práticas."
$ws.Range("C20").Value = "O aluno poderá optar por dois critérios de avaliação:Critério 1: NF = (P1+P2)/2; ouCritério 2: NF = (NOTA 1 + NOTA 2)/2Sendo P1 e P2 avaliações escritas e NOTA 1 e NOTA 2 obtidas em atividades desenvolvidas em aula, trabalhos e relatórios de aulas práticas."

$ws.Range("B21").Value = "Exame Final (EF) para alunos com Nota Final (NF) maior ou igual a 3,0 e menor do que 6,5 e frequência superior a 70%. Será considerado aprovado o aluno que tenha obtido Média Final (MF) igual ou maior do que 5,0, sendo MF = (NF+EF)/2."
$ws.Range("C21").Value = "Exame Final (EF) para alunos com Nota Final (NF) maior ou igual a 3,0 e menor do que 6,5 e frequência superior a 70%. Será considerado aprovado o aluno que tenha obtido Média Final (MF) igual ou maior do que 5,0, sendo MF = (NF+EF)/2."

$ws.Range("B22").Value = "Bibliografia básica:1. LEPSCH, I.F. 19 Lições de pedologia. São Paulo, Oficina do Texto. 456p. 2011. ISBN 978-85-7975-029-8.Bibliografia complementar:1. CAMARGO, O.A. de; MONIZ, A.C.; JORGE, J.A.; VALADARES, J.M.A.S. Métodos de analise química, mineralógica e física de solos do Instituto Agronômico de Campinas. Campinas, Instituto Agronômico, 2009. 77 p. (Boletim técnico, 106, Edição revista e atualizada).2. DIAS Jr., M.S. Compactação do solo. In: Tópicos em ciência do solo, v.1. NOVAIS, R.F.; ALVAREZ, V.H.; SCHAEFER, G.R. (Eds.). Viçosa: SBCS, 2000. p.55-94.3. EMBRAPA – EMPRESA BRASILEIRA DE PESQUISA AGROPECUÁRIA. Manual de análises químicas de solos, plantas e fertilizantes. SILVA, F. C. da (org.). EMBRAPA Comunicação para Transferência de Tecnologia. Brasília: EMBRAPA, Solos, 1999b. 370p.4. OLIVEIRA, J.B. Pedologia aplicada. 3a Edição. Piracicaba: Ed. FEALQ, 2008. 592p.5. REICHARDT, K.; TIMM, L.C. Solo, planta e atmosfera: conceitos, processos e aplicações. Barueri: SP: ed. Manole, 2004. 478p.6. SCHNEIDER, P.; GIASSON, E.; KLAMT, E. Classificação da aptidão agrícola das terras: um sistema alternativo. Porto Alegre: UFRGS, 2007. 72p."
$ws.Range("C22").Value = "Bibliografia básica:1. LEPSCH, I.F. 19 Lições de pedologia. São Paulo, Oficina do Texto. 456p. 2011. ISBN 978-85-7975-029-8.Bibliografia complementar:1. CAMARGO, O.A. de; MONIZ, A.C.; JORGE, J.A.; VALADARES, J.M.A.S. Métodos de analise química, mineralógica e física de solos do Instituto Agronômico de Campinas. Campinas, Instituto Agronômico, 2009. 77 p. (Boletim técnico, 106, Edição revista e atualizada).2. DIAS Jr., M.S. Compactação do solo. In: Tópicos em ciência do solo, v.1. NOVAIS, R.F.; ALVAREZ, V.H.; SCHAEFER, G.R. (Eds.). Viçosa: SBCS, 2000. p.55-94.3. EMBRAPA – EMPRESA BRASILEIRA DE PESQUISA AGROPECUÁRIA. Manual de análises químicas de solos, plantas e fertilizantes. SILVA, F. C. da (org.). EMBRAPA Comunicação para Transferência de Tecnologia. Brasília: EMBRAPA, Solos, 1999b. 370p.4. OLIVEIRA, J.B. Pedologia aplicada. 3a Edição. Piracicaba: Ed. FEALQ, 2008. 592p.5. REICHARDT, K.; TIMM, L.C. Solo, planta e atmosfera: conceitos, processos e aplicações. Barueri: SP: ed. Manole, 2004. 478p.6. SCHNEIDER, P.; GIASSON, E.; KLAMT, E. Classificação da aptidão agrícola das terras: um sistema alternativo. Porto Alegre: UFRGS, 2007. 72p."

